$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C8").Value = "00:04:57 (00:07:27 Â± 00:03:31)"
$ws.Range("E8").Value = "[2, 3, 5, 7, 11, 13, 17, 29, 31, 37, 41, 43, 47, 59, 61, 67, 71]"

$ws.Range("B10").Value = "0.649 (0.595 Â± 0.028)"
$ws.Range("C10").Value = "00:04:29 (00:04:30 Â± 00:00:02)"

$ws.Range("B16").Value = "0.666 (0.599 Â± 0.026)"
$ws.Range("C16").Value = "00:00:18 (00:00:19 Â± 00:00:00)"
